# Refresh the cryptocurrency price/volume table (cryptos.xlsx).
# Generated from the scheduled GitHub Actions scrape commit
# ("Updated cryptos list ... with GitHub Actions").
#
# All data cells on the sheet are stored as *text*, even when the text
# happens to look like a number (e.g. "1.001", "42.99"). Assigning such a
# string straight to Range.Value would make Excel auto-convert it to a
# numeric value, which would silently change the cell's stored type/shape.
# To avoid that we briefly force the cell to Text format, assign the
# value, then clear the format again (ClearFormats) so we don't leave a
# stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") values that parse as numbers --------------------
$numericLookingPrices = @{
    'D4' = '1.001'
    'D5' = '324.72'
    'D6' = '1.000'
    'D7' = '0.5094'
    'D8' = '0.4149'
    'D9' = '0.08748'
    'D10' = '1.134'
    'D11' = '42.99'
    'D14' = '6.571'
    'D15' = '7.462'
    'D17' = '94.23'
    'D18' = '0.00001115'
    'D19' = '0.06523'
    'D20' = '18.96'
    'D21' = '1.000'
    'D22' = '6.216'
    'D24' = '11.94'
    'D25' = '2.226'
    'D27' = '22.36'
    'D28' = '162.74'
    'D29' = '2.410'
    'D30' = '131.74'
    'D31' = '1.134'
    'D33' = '6.075'
    'D34' = '3.832'
    'D35' = '1.351'
    'D36' = '0.02522'
    'D37' = '5.431'
    'D38' = '0.06602'
    'D39' = '12.45'
    'D40' = '9.081'
    'D41' = '0.2196'
    'D42' = '0.6638'
    'D43' = '1.230'
    'D44' = '13.58'
    'D45' = '0.6167'
    'D46' = '2.188'
    'D47' = '3.662'
    'D48' = '1.264'
    'D49' = '124.30'
    'D50' = '80.31'
    'D51' = '0.06897'
}
foreach ($ref in $numericLookingPrices.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingPrices[$ref]
    $cell.ClearFormats()
}

# --- Remaining cells: Coin (B), Link (C), Price (D, non-numeric-looking)
# and Volume(1h) (E) -- plain text, safe to assign directly ---------------
$plainTextUpdates = @{
    'D2' = '30.405.05'
    'E2' = '  +1.23%  '
    'D3' = '2.010.64'
    'E3' = '  +5.03%  '
    'E4' = '  -0.05%  '
    'E5' = '  +1.43%  '
    'E7' = '  +1.12%  '
    'E8' = '  +3.05%  '
    'E9' = '  +6.04%  '
    'E10' = '  +2.42%  '
    'E11' = '  +2.28%  '
    'E12' = '  +3.32%  '
    'D13' = '1.996.36'
    'E13' = '  +3.73%  '
    'E14' = '  +2.32%  '
    'E15' = '  +2.48%  '
    'E16' = '  +0.05%  '
    'E17' = '  +2.27%  '
    'E18' = '  +1.68%  '
    'E19' = '  +0.29%  '
    'E20' = '  +4.61%  '
    'E22' = '  +4.57%  '
    'D23' = '30.446.72'
    'E23' = '  +1.25%  '
    'E24' = '  +5.68%  '
    'E25' = '  +1.11%  '
    'D26' = '2.231.69'
    'E26' = '  +4.22%  '
    'E27' = '  +0.00%  '
    'E28' = '  +0.75%  '
    'E29' = '  +6.06%  '
    'E30' = '  +2.27%  '
    'E31' = '  +0.22%  '
    'E32' = '  +1.27%  '
    'E34' = '  +1.27%  '
    'E35' = '  +12.49%  '
    'E36' = '  +3.14%  '
    'E37' = '  +1.62%  '
    'E38' = '  +2.78%  '
    'E39' = '  +9.17%  '
    'B40' = 'FraxShare'
    'C40' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E40' = '  +2.96%  '
    'B41' = 'Algorand'
    'C41' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E41' = '  +1.64%  '
    'E42' = '  +2.35%  '
    'E43' = '  +0.89%  '
    'E44' = '  +1.27%  '
    'E45' = '  +2.77%  '
    'E46' = '  -0.46%  '
    'E47' = '  +0.69%  '
    'E48' = '  +4.14%  '
    'E49' = '  +0.74%  '
    'E50' = '  +1.84%  '
    'E51' = '  +1.55%  '
}
foreach ($ref in $plainTextUpdates.Keys) {
    $ws.Range($ref).Value = $plainTextUpdates[$ref]
}
